$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing X10/Y10 cells on the existing last row (trade data for 2016-10-10)
$ws.Range("X10").Value = -0.010002000000000066
$ws.Range("Y10").Value = "Down"

# Copy number formats from row 10 down to row 11 for the date and percentage columns
# so the new cells reuse the existing styles (m/d/yyyy h:mm and 0.00%) instead of
# minting new ones.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("S10:T10").Copy()
$ws.Range("S11:T11").PasteSpecial(-4122)

# Append a new row of traded data (2016-10-11)
$ws.Range("A11").Value = 42654.882106481484
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Neutral"
$ws.Range("D11").Value = 24
$ws.Range("E11").Value = 18896
$ws.Range("F11").Value = 1072
$ws.Range("G11").Value = 64
$ws.Range("H11").Value = 34
$ws.Range("I11").Value = 82
$ws.Range("J11").Value = 17
$ws.Range("K11").Value = 28236
$ws.Range("L11").Value = 209
$ws.Range("M11").Value = 110
$ws.Range("N11").Value = 52
$ws.Range("O11").Value = 11
$ws.Range("P11").Value = "Named"
$ws.Range("Q11").Value = 17.089518681678967
$ws.Range("R11").Value = -24.44
$ws.Range("S11").Value = -0.1101
$ws.Range("T11").Value = -0.0419
$ws.Range("U11").Value = 6.47
$ws.Range("V11").Value = 1.88
$ws.Range("W11").Value = -2
